$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/string number format on cells whose new values look numeric,
# so Excel stores them as text (matching the inlineStr values in the target workbook)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.436.03'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '1.725.65'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '243.44'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4930'
$ws.Range("E7").Value = '  +2.61%  '
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("D9").Value = '0.06203'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").Value = '1.726.83'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = '0.06999'
$ws.Range("E11").Value = '  -2.52%  '
$ws.Range("D12").Value = '15.46'
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").Value = '4.551'
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = '77.50'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '0.9997'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '26.437.13'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '0.9995'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").Value = '0.000007171'
$ws.Range("E19").Value = '  +3.01%  '
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("D21").Value = '1.952.12'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '4.484'
$ws.Range("D23").Value = '8.592'
$ws.Range("E23").Value = '  -2.22%  '
$ws.Range("D24").Value = '5.169'
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("D25").Value = '138.01'
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").Value = '15.27'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '1.398'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").Value = '107.04'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").Value = '1.724'
$ws.Range("E29").Value = '  -2.89%  '
$ws.Range("D30").Value = '3.956'
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").Value = '0.08013'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '3.686'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").Value = '0.04522'
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.601'
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '0.9998'
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6286'
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.9363'
$ws.Range("E37").Value = '  +3.28%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '1.963'
$ws.Range("E38").Value = '  -4.27%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.390'
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '0.9997'
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01491'
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '99.54'
$ws.Range("E42").Value = '  -3.58%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.350'
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3857'
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.845'
$ws.Range("E45").Value = '  -2.82%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1168'
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05367'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '7.732'
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '30.22'
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.229'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '50.95'
$ws.Range("E51").Value = '  -0.62%  '
